$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.607758522033691
$ws.Range("B1").Value = 3.168262958526611
$ws.Range("C1").Value = 2.784998178482056
$ws.Range("D1").Value = 3.027395486831665
$ws.Range("E1").Value = 2.550008058547974
